$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Union): strip leading spaces from each quoted phrase in C2,
# and collapse the per-word cleaned list in E2 into comma-joined groups.
$ws.Range("C2").Value = '[''COSATU'', ''FOSATU'', ''NEHAWU'', ''SATAWU'', ''Allied'', ''AMCU'', ''Nuhhrccaw'', ''Denosa'', ''HOSPERSA'', ''NPSWU'', ''NUPSAW'', ''SADNU'', ''POPCRU'', ''NASUWU'', ''IMATU'', ''Union'', ''Shop stewart'', ''Shop steward'', ''Shopsteward'', ''Shopstewart'', ''Shop-steward'', ''Shop-stewart'']'
$ws.Range("E2").Value = '[''cosatu,fosatu,nehawu,satawu,allied,amcu,nuhhrccaw,denosa,hospersa,npswu,nupsaw,sadnu,popcru,nasuwu,imatu,union,shop'', ''stewart,shop'', ''steward,shopsteward,shopstewart,shop-steward,shop-stewart'']'

# Row 3 (Political Party)
$ws.Range("C3").Value = '[''ANC'', ''Congress'', ''DA'', ''EFF'', ''IFP'', ''NFP'', ''Cope'', ''UDM'', ''SACP'', ''ACDP'']'
$ws.Range("E3").Value = '[''anc,congress,da,eff,ifp,nfp,cope,udm,sacp,acdp'']'

# Row 4 (Civic org): also drops the trailing empty element after stripping.
$ws.Range("C4").Value = '[''association'', ''residents'', ''concerned'', ''ratepayers'', ''taxpayers'', ''NGO'', ''organisation'']'
$ws.Range("E4").Value = '[''association,residents,concerned,ratepayers,taxpayers,ngo,organisation'']'

# Row 5 (Church)
$ws.Range("C5").Value = '[''church'', ''congregation'', ''believers'', ''temple'', ''mosque'', ''synagoge'']'
$ws.Range("E5").Value = '[''church,congregation,believers,temple,mosque,synagoge'']'
